$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4430451.839000019
$ws.Range("C2").Value = 1841392.738000003
$ws.Range("D2").Value = 284335.0719999999
$ws.Range("E2").Value = 4488496.80300001
$ws.Range("F2").Value = 1987895.266000001
$ws.Range("G2").Value = 241172.399

$ws.Range("B3").Value = 4813409.601999978
$ws.Range("C3").Value = 1826769.926
$ws.Range("D3").Value = 268117.2559999998
$ws.Range("E3").Value = 4888585.296999977
$ws.Range("F3").Value = 1985415.666999999
$ws.Range("G3").Value = 232648.746

$ws.Range("B4").Value = 4798895.866999996
$ws.Range("C4").Value = 1825895.086999996
$ws.Range("D4").Value = 252090.1770000001
$ws.Range("E4").Value = 5012435.837000033
$ws.Range("F4").Value = 2104069.770999999
$ws.Range("G4").Value = 241855.4922120001

$ws.Range("B5").Value = 5459826.215999954
$ws.Range("C5").Value = 1829298.858000002
$ws.Range("D5").Value = 255064.4958420001
$ws.Range("E5").Value = 5402020.475999954
$ws.Range("F5").Value = 2306473.235000007
$ws.Range("G5").Value = 271051.334642

$ws.Range("B6").Value = 5078367.331999972
$ws.Range("C6").Value = 1838765.289000001
$ws.Range("D6").Value = 294418.201
$ws.Range("E6").Value = 6088796.679999981
$ws.Range("F6").Value = 2456039.330000006
$ws.Range("G6").Value = 297127.554

$ws.Range("B7").Value = 4534504.109000016
$ws.Range("C7").Value = 1958333.170000002
$ws.Range("D7").Value = 274656.1250000001
$ws.Range("E7").Value = 5760787.090999971
$ws.Range("F7").Value = 2727369.607000004
$ws.Range("G7").Value = 313127.383

$wb.Save()
